$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

# --- Paragraph: "how to ignore .class files in .gitingonre" ---
# Split into runs with gramStart/spellStart proofing-error marks.
$p = $d.Paragraphs(2)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t xml:space="preserve">how to ignore .class files </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>in .</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>gitingonre</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$rng.InsertXML($xml) | Out-Null

# --- Paragraph: "task2.2. checked ... and how to use them to find maximum,etc and what value switch statements accept" ---
# Split the "maximum,etc" run and add spellStart/spellEnd proofing marks, plus a trailing space run.
$p = $d.Paragraphs(4)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t>task2</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">.2. checked </w:t></w:r>' +
    '<w:hyperlink r:id="rId4" w:history="1">' +
        '<w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://www.tutorialspoint.com/java/switch_statement_in_java.htm</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>to how to create arrays</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">and how to use them to find </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>maximum,etc</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>and what value switch statements accept</w:t></w:r>' +
    '</w:p>'
$rng.InsertXML($xml) | Out-Null

# --- Paragraph: weeklyPay error message ---
# Split "weeklyPay" and "weeklypay" off into their own runs with spellStart/spellEnd marks.
$p = $d.Paragraphs(5)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t>I got this error saying “</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">The local variable </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>weeklyPay</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> may not have been initialized</w:t></w:r>' +
    '<w:r><w:t>” I solved it by assigning a value to it. “</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>weeklypay</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>=0;”</w:t></w:r>' +
    '</w:p>'
$rng.InsertXML($xml) | Out-Null

# --- Insert one extra empty paragraph right after the weeklyPay paragraph ---
$p6 = $d.Paragraphs(6)
$insertPoint = $d.Range($p6.Range.Start, $p6.Range.Start)
$insertPoint.InsertXML('<w:p ' + $wNs + '/>') | Out-Null
